$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowA = $ws.Range("B39:AD39").Value2
$rowB = $ws.Range("B40:AD40").Value2
$ws.Range("B39:AD39").Value2 = $rowB
$ws.Range("B40:AD40").Value2 = $rowA

$rowA = $ws.Range("B45:AD45").Value2
$rowB = $ws.Range("B46:AD46").Value2
$ws.Range("B45:AD45").Value2 = $rowB
$ws.Range("B46:AD46").Value2 = $rowA

$rowA = $ws.Range("B54:AD54").Value2
$rowB = $ws.Range("B55:AD55").Value2
$ws.Range("B54:AD54").Value2 = $rowB
$ws.Range("B55:AD55").Value2 = $rowA

$rowA = $ws.Range("B62:AD62").Value2
$rowB = $ws.Range("B63:AD63").Value2
$ws.Range("B62:AD62").Value2 = $rowB
$ws.Range("B63:AD63").Value2 = $rowA

$rowA = $ws.Range("B81:AD81").Value2
$rowB = $ws.Range("B82:AD82").Value2
$ws.Range("B81:AD81").Value2 = $rowB
$ws.Range("B82:AD82").Value2 = $rowA

$rowA = $ws.Range("B95:AD95").Value2
$rowB = $ws.Range("B96:AD96").Value2
$ws.Range("B95:AD95").Value2 = $rowB
$ws.Range("B96:AD96").Value2 = $rowA

$rowA = $ws.Range("B100:AD100").Value2
$rowB = $ws.Range("B101:AD101").Value2
$ws.Range("B100:AD100").Value2 = $rowB
$ws.Range("B101:AD101").Value2 = $rowA

$rowA = $ws.Range("B117:AD117").Value2
$rowB = $ws.Range("B119:AD119").Value2
$ws.Range("B117:AD117").Value2 = $rowB
$ws.Range("B119:AD119").Value2 = $rowA

$rowA = $ws.Range("B126:AD126").Value2
$rowB = $ws.Range("B127:AD127").Value2
$ws.Range("B126:AD126").Value2 = $rowB
$ws.Range("B127:AD127").Value2 = $rowA

$rowA = $ws.Range("B140:AD140").Value2
$rowB = $ws.Range("B141:AD141").Value2
$ws.Range("B140:AD140").Value2 = $rowB
$ws.Range("B141:AD141").Value2 = $rowA

$rowA = $ws.Range("B148:AD148").Value2
$rowB = $ws.Range("B150:AD150").Value2
$ws.Range("B148:AD148").Value2 = $rowB
$ws.Range("B150:AD150").Value2 = $rowA

$rowA = $ws.Range("B152:AD152").Value2
$rowB = $ws.Range("B153:AD153").Value2
$ws.Range("B152:AD152").Value2 = $rowB
$ws.Range("B153:AD153").Value2 = $rowA

$rowA = $ws.Range("B161:AD161").Value2
$rowB = $ws.Range("B162:AD162").Value2
$ws.Range("B161:AD161").Value2 = $rowB
$ws.Range("B162:AD162").Value2 = $rowA

$rowA = $ws.Range("B187:AD187").Value2
$rowB = $ws.Range("B188:AD188").Value2
$ws.Range("B187:AD187").Value2 = $rowB
$ws.Range("B188:AD188").Value2 = $rowA

$rowA = $ws.Range("B217:AD217").Value2
$rowB = $ws.Range("B218:AD218").Value2
$ws.Range("B217:AD217").Value2 = $rowB
$ws.Range("B218:AD218").Value2 = $rowA
